$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")
$r = $ws.Range("Z10:AB12")
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(7).Weight = -4138
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(10).Weight = 2
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(8).Weight = 2
$r.Borders.Item(9).LineStyle = 1
$r.Borders.Item(9).Weight = 2
$r.Borders.Item(11).LineStyle = 1
$r.Borders.Item(11).Weight = 2
$r.Borders.Item(12).LineStyle = 1
$r.Borders.Item(12).Weight = 2
Write-Host "done"
